# Add data for 2024-05-20
# Updates the running 2024 (column K) totals — and a few backfilled
# 2015-2022 corrections in column B/I — across the citywide, by-neighborhood
# summary sheets and each individual neighborhood sheet.

$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("Citywide Totals")
$ws.Range("K2").Value = 2842
$ws.Range("K3").Value = 2771
$ws.Range("B4").Value = 1699
$ws.Range("I4").Value = 1791
$ws.Range("K4").Value = 579
$ws.Range("K6").Value = 3387
$ws.Range("B7").Value = 23332
$ws.Range("I7").Value = 26244
$ws.Range("K7").Value = 9759

$ws = $wb.Worksheets.Item("Logan Square")
$ws.Range("K2").Value = 30
$ws.Range("K7").Value = 140

$ws = $wb.Worksheets.Item("Austin")
$ws.Range("K2").Value = 186
$ws.Range("K3").Value = 194
$ws.Range("K4").Value = 35
$ws.Range("K6").Value = 213
$ws.Range("K7").Value = 642

$ws = $wb.Worksheets.Item("South Chicago")
$ws.Range("K3").Value = 66
$ws.Range("K6").Value = 51
$ws.Range("K7").Value = 213

$ws = $wb.Worksheets.Item("Garfield Park")
$ws.Range("K2").Value = 107
$ws.Range("K3").Value = 140
$ws.Range("K4").Value = 21
$ws.Range("K6").Value = 109
$ws.Range("K7").Value = 386

$ws = $wb.Worksheets.Item("West Pullman")
$ws.Range("K2").Value = 50
$ws.Range("K3").Value = 55
$ws.Range("K7").Value = 160

$ws = $wb.Worksheets.Item("Grand Crossing")
$ws.Range("K2").Value = 84
$ws.Range("K3").Value = 110
$ws.Range("K6").Value = 100
$ws.Range("K7").Value = 321

$ws = $wb.Worksheets.Item("New City")
$ws.Range("K2").Value = 69
$ws.Range("K6").Value = 97
$ws.Range("K7").Value = 230

$ws = $wb.Worksheets.Item("Woodlawn")
$ws.Range("K2").Value = 53
$ws.Range("K3").Value = 65
$ws.Range("K7").Value = 176

$ws = $wb.Worksheets.Item("Fuller Park")
$ws.Range("K3").Value = 11
$ws.Range("K7").Value = 29

$ws = $wb.Worksheets.Item("By Neighborhood")
$ws.Range("K2").Value = 75
$ws.Range("K7").Value = 286
$ws.Range("K8").Value = 642
$ws.Range("K11").Value = 206
$ws.Range("K19").Value = 290
$ws.Range("K20").Value = 225
$ws.Range("K23").Value = 85
$ws.Range("K29").Value = 505
$ws.Range("K30").Value = 29
$ws.Range("K33").Value = 386
$ws.Range("K36").Value = 113
$ws.Range("K37").Value = 321
$ws.Range("K42").Value = 342
$ws.Range("K43").Value = 86
$ws.Range("K44").Value = 93
$ws.Range("K47").Value = 51
$ws.Range("K48").Value = 117
$ws.Range("K52").Value = 272
$ws.Range("K53").Value = 140
$ws.Range("K54").Value = 183
$ws.Range("K55").Value = 105
$ws.Range("K60").Value = 62
$ws.Range("B63").Value = 404
$ws.Range("I63").Value = 204
$ws.Range("K63").Value = 44
$ws.Range("K64").Value = 60
$ws.Range("K65").Value = 230
$ws.Range("K67").Value = 383
$ws.Range("K73").Value = 93
$ws.Range("K76").Value = 150
$ws.Range("K77").Value = 70
$ws.Range("K82").Value = 11
$ws.Range("K83").Value = 213
$ws.Range("K85").Value = 466
$ws.Range("K86").Value = 63
$ws.Range("K88").Value = 111
$ws.Range("K89").Value = 131
$ws.Range("K91").Value = 92
$ws.Range("K94").Value = 116
$ws.Range("K95").Value = 160
$ws.Range("K99").Value = 176
$ws.Range("B101").Value = 23332
$ws.Range("I101").Value = 26244
$ws.Range("K101").Value = 9759

$ws = $wb.Worksheets.Item("North Lawndale")
$ws.Range("K2").Value = 121
$ws.Range("K3").Value = 122
$ws.Range("K6").Value = 111
$ws.Range("K7").Value = 383

$ws = $wb.Worksheets.Item("Loop")
$ws.Range("K6").Value = 79
$ws.Range("K7").Value = 183

$ws = $wb.Worksheets.Item("Englewood")
$ws.Range("K2").Value = 138
$ws.Range("K3").Value = 168
$ws.Range("K6").Value = 159
$ws.Range("K7").Value = 505

$ws = $wb.Worksheets.Item("Lake View")
$ws.Range("K6").Value = 61
$ws.Range("K7").Value = 117

$ws = $wb.Worksheets.Item("Chatham")
$ws.Range("K2").Value = 93
$ws.Range("K3").Value = 78
$ws.Range("K6").Value = 95
$ws.Range("K7").Value = 290

$ws = $wb.Worksheets.Item("Irving Park")
$ws.Range("K6").Value = 44
$ws.Range("K7").Value = 93

$ws = $wb.Worksheets.Item("River North")
$ws.Range("K3").Value = 27
$ws.Range("K6").Value = 89
$ws.Range("K7").Value = 150

$ws = $wb.Worksheets.Item("Humboldt Park")
$ws.Range("K6").Value = 132
$ws.Range("K7").Value = 342

$ws = $wb.Worksheets.Item("Lower West Side")
$ws.Range("K2").Value = 37
$ws.Range("K7").Value = 105

$ws = $wb.Worksheets.Item("Douglas")
$ws.Range("K3").Value = 26
$ws.Range("K6").Value = 24
$ws.Range("K7").Value = 85

$ws = $wb.Worksheets.Item("Washington Park")
$ws.Range("K2").Value = 26
$ws.Range("K3").Value = 41
$ws.Range("K7").Value = 92

$ws = $wb.Worksheets.Item("Near South Side")
$ws.Range("K2").Value = 14
$ws.Range("K7").Value = 60

$ws = $wb.Worksheets.Item("Chicago Lawn")
$ws.Range("K3").Value = 62
$ws.Range("K6").Value = 79
$ws.Range("K7").Value = 225

$ws = $wb.Worksheets.Item("Grand Boulevard")
$ws.Range("K2").Value = 44
$ws.Range("K6").Value = 26
$ws.Range("K7").Value = 113

$ws = $wb.Worksheets.Item("Auburn Gresham")
$ws.Range("K2").Value = 98
$ws.Range("K4").Value = 12
$ws.Range("K7").Value = 286

$ws = $wb.Worksheets.Item("West Loop")
$ws.Range("K3").Value = 23
$ws.Range("K7").Value = 116

$ws = $wb.Worksheets.Item("Kenwood")
$ws.Range("K3").Value = 18
$ws.Range("K7").Value = 51

$ws = $wb.Worksheets.Item("Belmont Cragin")
$ws.Range("K2").Value = 62
$ws.Range("K3").Value = 54
$ws.Range("K7").Value = 206

$ws = $wb.Worksheets.Item("Portage Park")
$ws.Range("K2").Value = 27
$ws.Range("K7").Value = 93

$ws = $wb.Worksheets.Item("Albany Park")
$ws.Range("K3").Value = 21
$ws.Range("K7").Value = 75

$ws = $wb.Worksheets.Item("United Center")
$ws.Range("K2").Value = 26
$ws.Range("K7").Value = 111

$ws = $wb.Worksheets.Item("Uptown")
$ws.Range("K2").Value = 31
$ws.Range("K3").Value = 42
$ws.Range("K4").Value = 20
$ws.Range("K7").Value = 131

$ws = $wb.Worksheets.Item("Streeterville")
$ws.Range("K4").Value = 23
$ws.Range("K7").Value = 63

$ws = $wb.Worksheets.Item("Morgan Park")
$ws.Range("K2").Value = 18
$ws.Range("K7").Value = 62

$ws = $wb.Worksheets.Item("Hyde Park")
$ws.Range("K3").Value = 24
$ws.Range("K6").Value = 37
$ws.Range("K7").Value = 86

$ws = $wb.Worksheets.Item("South Shore")
$ws.Range("K3").Value = 161
$ws.Range("K7").Value = 466

$ws = $wb.Worksheets.Item("Sheffield & DePaul")
$ws.Range("K4").Value = 1
$ws.Range("K6").Value = 11

$ws = $wb.Worksheets.Item("Riverdale")
$ws.Range("K3").Value = 25
$ws.Range("K7").Value = 70

$ws = $wb.Worksheets.Item("Little Village")
$ws.Range("K2").Value = 74
$ws.Range("K3").Value = 66
$ws.Range("K7").Value = 272
